$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns before column D (old D:K shift right to F:M)
# ---------------------------------------------------------------------------
$ws.Range("D:E").Insert()

# Carry over the number formatting (date format in row 7/38/80, number format
# elsewhere) from the old column D (now column F) into the two new columns.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Populate the new columns D (newest quarter) and E (next quarter) with
#    the new financial data.
# ---------------------------------------------------------------------------
$rows  = @(7,8,9,10,12,13,14,15,17,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,38,41,42,43,44,45,46,47,48,49,50,51,52,53,54,57,58,59,60,61,62,63,64,65,66,68,69,70,71,72,73,74,75,76,77,80,81,83,84,85,86,87,88,89,91,92,93,94,96,97,98,99,100,101,102)
$dvals = @(43465,38500,25800,12700,"NA",0,1200,1500,38500,0,0,1500,2900,-2900,500,0,-3500,-3500,0,"NA",0,0,0,-3500,0,-3500,43465,3500,0,19000,0,2500,25000,0,800,104700,0,0,400,0,130900,2900,3200,7400,13500,68800,1200,0,0,0,83500,0,0,27600,0,-26500,0,0,0,19900,0,43465,-3500,1500,0,0,0,0,0,800,0,0,0,0,0,0,0,0,-500,0,300)
$evals = @(43373,39900,25100,14800,"NA",0,800,1500,38900,1000,0,2500,3100,-2100,-700,0,-1400,-1400,0,400,0,0,0,-1000,0,-1000,43373,3200,0,20800,0,2300,26200,0,900,106100,0,0,400,0,133600,2500,2400,8200,13100,69900,700,0,0,0,83700,0,0,28800,0,-23000,0,0,0,21100,0,43373,-1000,1500,0,0,0,0,0,800,-100,0,0,-100,0,0,0,0,-100,0,600)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
    $ws.Cells.Item($r, 5).Value = $evals[$i]
}

# ---------------------------------------------------------------------------
# 3. A handful of cells in the shifted range received genuinely new figures
#    (not merely the shifted-over value) - correct those explicitly.
# ---------------------------------------------------------------------------
$ws.Range("I59").Value = 13800
$ws.Range("I60").Value = 19500
$ws.Range("I61").Value = 67600
$ws.Range("I62").Value = 1700

# ---------------------------------------------------------------------------
# 4. Column widths - keep them close to their previous look after the new
#    columns were inserted.
# ---------------------------------------------------------------------------
$ws.Columns("D:E").ColumnWidth = 14.66
$ws.Columns("F").ColumnWidth = 14.44
$ws.Columns("G").ColumnWidth = 14.89
$ws.Columns("H:I").ColumnWidth = 14.66
$ws.Columns("J").ColumnWidth = 14.44
$ws.Columns("K").ColumnWidth = 14.89
$ws.Columns("L:M").ColumnWidth = 14.66

Write-Host "Edit complete"
